$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New log entries for 17.10.2022 documenting ProjectLibre research/tutorials,
# added right after the existing "Methodologies" row (row 32). Copy the
# formatting of the row above (date-column shading + time-column format)
# before filling in the new values.
$ws.Range("A32:F32").Copy()
$ws.Range("A33:F34").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A33").Value = "17.10.2022"
$ws.Range("B33").Value = 0.83333333333333337
$ws.Range("C33").Value = "ProjectLibre"
$ws.Range("D33").Value = "Tutorial"
$ws.Range("E33").Value = 20
$ws.Range("F33").Value = "Tutorials on How to Use ProjectLibre Software"

$ws.Range("A34").Value = "17.10.2022"
$ws.Range("B34").Value = 0.84722222222222221
$ws.Range("C34").Value = "ProjectLibre"
$ws.Range("D34").Value = "Documentation"
$ws.Range("E34").Value = 100
$ws.Range("F34").Value = "Create Work Breakdown Structure and Gantt Chart with Project Libre"

# Move the saved selection to match where the author ended up after the edit.
[void]$ws.Range("F35").Select()
